$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. header { ... } : "padding: 1rem;" -> "padding: 1rem 2rem;"
#    (text is not unique document-wide, so anchor on the preceding line)
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "justify-content: center;^p      padding: 1rem;", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "justify-content: center;^p      padding: 1rem 2rem;", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. header h1 { ... } : "font-size: 2rem;" -> "font-size: 2.2rem;"
#    (this text is unique in the document)
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "font-size: 2rem;", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "font-size: 2.2rem;", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. header p { ... } : "margin: 0;" -> "margin: 0.2rem 0 0;"
#    (text is not unique document-wide, so anchor on the preceding line)
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "header p {^p      margin: 0;", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "header p {^p      margin: 0.2rem 0 0;", 2) | Out-Null

# ---------------------------------------------------------------------
# 4. Insert a new "@media screen and (max-width: 600px) { ... }" block
#    right after the existing "footer { ... }" block, before "</style>".
#    Blank CSS lines are inserted as a unique placeholder token first and
#    then cleared out to a true empty paragraph (<w:p/>) afterwards, which
#    is what Word produces for a manually-typed blank line.
# ---------------------------------------------------------------------
$oldTail = "    footer {^p      background: #003366;^p      color: white;^p      text-align: center;^p      padding: 1rem;^p    }^p  </style>"
$newTail = "    footer {^p      background: #003366;^p      color: white;^p      text-align: center;^p      padding: 1rem;^p    }^p@@BLANK1@@^p    @media screen and (max-width: 600px) {^p      header img.logo {^p        height: 40px;^p        left: 10px;^p      }^p@@BLANK2@@^p      header h1 {^p        font-size: 1.5rem;^p      }^p@@BLANK3@@^p      header p {^p        font-size: 0.9rem;^p      }^p@@BLANK4@@^p      section {^p        padding: 1rem;^p      }^p    }^p  </style>"

$d.Content.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2) | Out-Null

foreach ($ph in @("@@BLANK1@@", "@@BLANK2@@", "@@BLANK3@@", "@@BLANK4@@")) {
    $r = $d.Content
    $found = $r.Find.Execute($ph, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $r.Text = ""
    }
}

# ---------------------------------------------------------------------
# 5. "<h2>Contact Us</h2>" paragraph gains a right tab stop at 8280 twips
#    (414 pt) plus a tab character, and the document's lone "_GoBack"
#    bookmark (previously sitting at the very end of the document, after
#    "</html>") moves to sit right after that new tab.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("    <h2>Contact Us</h2>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $r.Paragraphs(1)
$para.TabStops.Add(414, 0) | Out-Null

# Insert the tab character, followed by a throw-away marker so we have a
# non-boundary insertion point to anchor the relocated bookmark on.
$d.Content.Find.Execute(
    "    <h2>Contact Us</h2>", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "    <h2>Contact Us</h2>^t@@GOBACK@@", 2) | Out-Null

$bmRange = $d.Content
$bmRange.Find.Execute("@@GOBACK@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$d.Content.Find.Execute("@@GOBACK@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------
# 6. Footer copyright year: 2025 -> 2022
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "&copy; 2025 The Learning Curve.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "&copy; 2022 The Learning Curve.", 2) | Out-Null
